$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 data (09-09-2021)
# Force text for the date-looking value so it is stored as a string
# (matching existing cells A2/A3 which are plain shared-string text, no style).
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "09-09-2021"
$ws.Cells.Item(4, 1).Style = "Normal"

$ws.Cells.Item(4, 2).Value = 320000
$ws.Cells.Item(4, 3).Value = 535000
$ws.Cells.Item(4, 4).Value = 320000
$ws.Cells.Item(4, 5).Value = 315000
$ws.Cells.Item(4, 6).Value = 5000
$ws.Cells.Item(4, 7).Value = 4.15
